$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("heating_costs")

# --- Copy formatting of the last existing data row (9) down into the two
#     new rows (10, 11) so fills/fonts/number-formats match the existing
#     table styling exactly, then overwrite with the new row's content. ---
$ws.Range("A9:S9").Copy()
$ws.Range("A10:S11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New text values are introduced in this precise order so the shared-string
# table is built up the same way it was authored: technology, then the
# notes column, then each row's distinct efficiency label, then the brand.
$ws.Cells.Item(10, 2).Value = "Electric ASHP - Ducted"
$ws.Cells.Item(10, 19).Value = "Data Year 2025"
$ws.Cells.Item(10, 4).Value = "SEER 15, 8.8 HSPF"
$ws.Cells.Item(11, 4).Value = "SEER 16, 8.8 HSPF"
$ws.Cells.Item(10, 18).Value = "Trane"

# Row 10: Install / Electric ASHP - Ducted / Electricity / SEER 15, 8.8 HSPF
$ws.Cells.Item(10, 1).Value = "Install"
$ws.Cells.Item(10, 3).Value = "Electricity"
$ws.Cells.Item(10, 5).Value = 2025
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1
$ws.Cells.Item(10, 8).Value = 6654.7500000000018
$ws.Cells.Item(10, 9).Value = 8151.0350000000008
$ws.Cells.Item(10, 10).Value = 9647.32
$ws.Cells.Item(10, 11).Value = 87.852499999999964
$ws.Cells.Item(10, 12).Value = 95.945416666666631
$ws.Cells.Item(10, 13).Value = 104.03833333333331
$ws.Cells.Item(10, 14).Value = 0
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Formula = "=`$O10+(`$O10-`$N10)"
$ws.Cells.Item(10, 17).Value = 25

# Row 11: Install / Electric ASHP - Ducted / Electricity / SEER 16, 8.8 HSPF
$ws.Cells.Item(11, 1).Value = "Install"
$ws.Cells.Item(11, 2).Value = "Electric ASHP - Ducted"
$ws.Cells.Item(11, 3).Value = "Electricity"
$ws.Cells.Item(11, 5).Value = 2025
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1
$ws.Cells.Item(11, 8).Value = 7864.9399999999987
$ws.Cells.Item(11, 9).Value = 8664.32
$ws.Cells.Item(11, 10).Value = 9463.6999999999989
$ws.Cells.Item(11, 11).Value = 93.546666666666681
$ws.Cells.Item(11, 12).Value = 102.6179166666667
$ws.Cells.Item(11, 13).Value = 111.68916666666671
$ws.Cells.Item(11, 14).Value = 0
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Formula = "=`$O11+(`$O11-`$N11)"
$ws.Cells.Item(11, 17).Value = 25
$ws.Cells.Item(11, 18).Value = "Trane"
$ws.Cells.Item(11, 19).Value = "Data Year 2025"

# --- Selection / active-sheet bookkeeping: the author ended the session
#     with heating_costs active (instead of cpi) and the cursor on Q12. ---
$ws.Activate()
$ws.Range("Q12").Select()
